$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 632 (shifts rows 632..673 down to 633..674)
$ws.Range("A632").EntireRow.Insert()

# Fill in the new row's data
$ws.Range("A632").NumberFormat = "@"
$ws.Range("A632").Value = "2026/01/12"
$ws.Range("A632").ClearFormats()
$ws.Range("B632").Value = "月"
$ws.Range("C632").Value = 4
$ws.Range("D632").Value = 139
